# Adds 20 new User Story rows (162-181) for the new 'Transportations' domain
# within the 'Social and Urban Studies' domain cluster, matching the commit:
# Add USs related to cluster domain 'Social and Urban Studies' [Transportations]

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$domainCluster = 'Social and Urban Studies'
$domain = 'Transportations'
$promptUsed = 'Domain_FSPrompt'

$tasks = @(
    'adversarial learning',
    'cnn',
    'conversational agent',
    'decision tree',
    'document classification',
    'entity extraction',
    'feature selection',
    'imbalanced dataset',
    'k-nearest neighbor',
    'keyword extraction',
    'multi-label classification',
    'neural network',
    'random forest',
    'semantic similarity',
    'sentiment analysis',
    'speech to text',
    'text categorization',
    'unsupervised clustering',
    'voice recognition',
    'word embedding'
)

$stories = @(
    'As a logistics manager, I want to employ adversarial learning models to identify fraudulent activities in shipment tracking systems, so that I can prevent cargo theft and ensure the integrity of supply chain operations.',
    'As a railway operator, I want to employ CNN architectures to analyze track condition images captured by drones, so that I can detect defects early, schedule maintenance proactively, and improve overall safety and reliability.',
    'As an airport information desk operator, I want to implement a conversational agent to help travelers navigate through the airport, find amenities, and receive real-time updates on flight statuses, so that I can enhance customer satisfaction and streamline operations.',
    'As a traffic engineer, I want to use decision tree models to analyze historical traffic data and identify key factors contributing to congestion at specific intersections, so that I can optimize traffic signal timings and improve traffic flow.',
    'As a safety compliance officer for a transportation company, I want to implement document classification models to categorize and prioritize incident reports and safety inspection records, so that I can identify trends, assess risks, and improve safety protocols.',
    'As an airline operations manager, I want to deploy entity extraction models to automatically extract flight details including departure times, arrival times, flight numbers, and destinations from flight schedules and booking confirmations, so that I can monitor flight statuses and manage operations more effectively.',
    'As a logistics planner, I want to use feature selection techniques to identify key variables in shipment data that impact delivery times, such as distance, mode of transport, and traffic conditions, so that I can improve delivery efficiency and customer satisfaction.',
    'As a railway safety officer, I want to address the imbalanced dataset in incident reports to improve the accuracy of predicting rare but critical safety incidents, such as derailments or collisions, so that I can implement targeted safety measures and reduce risks.',
    'As a public transport planner, I want to utilize k-NN algorithms to recommend optimal transit routes to commuters based on their historical travel patterns and current traffic conditions, so that I can enhance commuter satisfaction and promote the use of public transportation.',
    'As a transportation planner, I want to apply keyword extraction methods to analyze public feedback and comments on transportation services to identify key concerns and suggestions, so that I can improve service quality and customer satisfaction.',
    'As a cargo logistics manager, I want to use multi-label classification to predict the types of goods being transported based on shipment descriptions, so that I can ensure proper handling, storage, and compliance with transport regulations.',
    'As an autonomous vehicle developer, I want to train neural networks to recognize and classify different types of road obstacles and hazards from sensor data, such as pedestrians, cyclists, and road debris, so that I can improve the vehicle''s ability to navigate safely.',
    'As a logistics manager, I want to deploy a random forest algorithm to classify delivery routes based on factors such as distance, road conditions, and delivery time windows, so that I can optimize route planning and ensure timely deliveries.',
    'As a traffic management system developer, I want to implement semantic similarity algorithms to identify similar traffic incident reports and prioritize response efforts based on the severity and impact on traffic flow, so that I can improve incident management efficiency.',
    'As a public transport planner, I want to use sentiment analysis to analyze commuter sentiment towards different routes and services based on social media posts and surveys, so that I can prioritize enhancements and adjustments in public transport offerings.',
    'As a railway station announcer, I want to deploy speech-to-text systems to transcribe train arrival and departure announcements in real-time, so that I can provide accurate and timely information to passengers, especially in noisy or crowded environments.',
    'As an airport operations manager, I want to utilize text categorization algorithms to categorize passenger feedback received through surveys and comment cards into categories such as facilities, staff interactions, and amenities, so that I can identify areas for improvement and enhance passenger experience.',
    'As a traffic management engineer, I want to use unsupervised clustering algorithms to group traffic sensors based on similar traffic patterns and congestion levels, so that I can identify critical areas for traffic flow optimization and signal timing adjustments.',
    'As a railway station announcer, I want to deploy voice recognition software to automatically transcribe and broadcast train departure and arrival announcements, so that I can ensure accurate and timely information for passengers.',
    'As a logistics planner, I want to use word embedding models to analyze customer feedback and identify key themes and sentiments related to delivery services, so that I can prioritize improvements and enhance customer satisfaction.'
)

$startRow = 162
for ($i = 0; $i -lt $tasks.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $domainCluster
    $ws.Cells.Item($r, 2).Value = 7
    $ws.Cells.Item($r, 3).Value = $domain
    $ws.Cells.Item($r, 4).Value = $tasks[$i]
    $ws.Cells.Item($r, 5).Value = $stories[$i]
    $ws.Cells.Item($r, 6).Value = $promptUsed
}

Write-Host 'Rows 162-181 populated'
